$wb = $excel.ActiveWorkbook

# --- Sheet1: add two new student rows, shorten B3, move selection ---
$ws1 = $wb.Worksheets.Item(1)

$ws1.Range("B3").Value = "ABCV"

$ws1.Range("A9").Value = "A001"
$ws1.Range("B9").Value = "ABC"
$ws1.Range("A10").Value = "A002"
$ws1.Range("B10").Value = "ACBD"

$ws1.Range("B3").Select()

# --- Sheet2: rename to the new test date, update/trim the score rows ---
$ws2 = $wb.Worksheets.Item(2)
$ws2.Name = "test 2021-09-22"

$ws2.Range("A1").Value = "test 2021-09-22"

$ws2.Range("A3").Value = "C020"
$ws2.Range("B3").Value = 48

$ws2.Range("A4").Value = "A021"
$ws2.Range("B4").Value = 44

$ws2.Range("A5:B8").ClearContents()
